# Auto-generated edit script: appends two new data rows (184, 185)
# for date 2025-12-01 to Sheet1, matching the existing rows for the
# "si charging volume" stations (si=四方坪站 / 高岭站).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 184: 2025-12-01 / 四方坪站充电量(kw)
$ws.Range("A184").Value = "2025-12-01"
$ws.Range("B184").Value = "四方坪站充电量(kw)"
$ws.Range("C184").Value = 579.11
$ws.Range("D184").Value = 917.38099999999997
$ws.Range("E184").Value = 636.51100000000008
$ws.Range("F184").Value = 284.62
$ws.Range("G184").Value = 242.64399999999998
$ws.Range("H184").Value = 641.32799999999997
$ws.Range("I184").Value = 335.28099999999995
$ws.Range("J184").Value = 172.67099999999999
$ws.Range("K184").Value = 93.634
$ws.Range("L184").Value = 309.59699999999998
$ws.Range("M184").Value = 81.89700000000002
$ws.Range("N184").Value = 195.886
$ws.Range("O184").Value = 493.92199999999997
$ws.Range("P184").Value = 1124.3910000000001
$ws.Range("Q184").Value = 393.84899999999999
$ws.Range("R184").Value = 326.8
$ws.Range("S184").Value = 358.39599999999996
$ws.Range("T184").Value = 305.74900000000002
$ws.Range("U184").Value = 109.434
$ws.Range("V184").Value = 42.23
$ws.Range("W184").Value = 57.8
$ws.Range("X184").Value = 105.85
$ws.Range("Y184").Value = 62.080000000000005
$ws.Range("Z184").Value = 66.62

# Row 185: 2025-12-01 / 高岭站充电量(kw)
$ws.Range("A185").Value = "2025-12-01"
$ws.Range("B185").Value = "高岭站充电量(kw)"
$ws.Range("C185").Value = 422.26799999999992
$ws.Range("D185").Value = 351.38799999999998
$ws.Range("E185").Value = 167.74799999999999
$ws.Range("F185").Value = 97.085999999999999
$ws.Range("G185").Value = 85.59
$ws.Range("H185").Value = 224.833
$ws.Range("I185").Value = 157.542
$ws.Range("J185").Value = 64.754999999999995
$ws.Range("K185").Value = 275.56299999999999
$ws.Range("L185").Value = 174.89999999999998
$ws.Range("M185").Value = 210.322
$ws.Range("N185").Value = 144.643
$ws.Range("O185").Value = 499.19900000000001
$ws.Range("P185").Value = 332.07500000000005
$ws.Range("Q185").Value = 189.904
$ws.Range("R185").Value = 225.28299999999999
$ws.Range("S185").Value = 415.20399999999995
$ws.Range("T185").Value = 21.422000000000004
$ws.Range("U185").Value = 0
$ws.Range("V185").Value = 54.030999999999999
$ws.Range("W185").Value = 105.36699999999999
$ws.Range("X185").Value = 0
$ws.Range("Y185").Value = 17.030999999999999
$ws.Range("Z185").Value = 0

# Restore the selection to match the new bottom of the data range
$ws.Range("D189").Select()
